$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellAddr, $val) {
    $r = $ws.Range($cellAddr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.ClearFormats()
}

Set-TextValue "D2" "26.840.90"
Set-TextValue "E2" "  +0.23%  "
Set-TextValue "D3" "1.642.31"
Set-TextValue "E3" "  -0.11%  "
Set-TextValue "E4" "  -0.21%  "
Set-TextValue "D5" "218.54"
Set-TextValue "E5" "  +0.88%  "
Set-TextValue "D6" "0.500"
Set-TextValue "E6" "  +0.06%  "
Set-TextValue "E7" "  -0.24%  "
Set-TextValue "E8" "  -0.12%  "
Set-TextValue "D9" "0.0622"
Set-TextValue "E9" "  -1.03%  "
Set-TextValue "D10" "19.25"
Set-TextValue "E10" "  +0.49%  "
Set-TextValue "D11" "0.0848"
Set-TextValue "E11" "  +0.94%  "
Set-TextValue "D12" "1.871.58"
Set-TextValue "E12" "  -0.06%  "
Set-TextValue "D13" "1.644.58"
Set-TextValue "E13" "  -0.08%  "
Set-TextValue "E14" "  -0.02%  "
Set-TextValue "E15" "  -0.34%  "
Set-TextValue "D16" "65.42"
Set-TextValue "E16" "  +1.75%  "
Set-TextValue "D17" "26.848.82"
Set-TextValue "E17" "  +0.21%  "
Set-TextValue "D18" "0.0₃0736"
Set-TextValue "E18" "  -0.38%  "
Set-TextValue "D19" "215.65"
Set-TextValue "E19" "  +0.91%  "
Set-TextValue "E20" "  -0.25%  "
Set-TextValue "E21" "  +6.20%  "
Set-TextValue "E22" "  +0.18%  "
Set-TextValue "E23" "  -0.14%  "
Set-TextValue "E24" "  -1.65%  "
Set-TextValue "D25" "147.83"
Set-TextValue "E25" "  +1.91%  "
Set-TextValue "E26" "  -0.23%  "
Set-TextValue "E27" "  +0.10%  "
Set-TextValue "E28" "  +1.02%  "
Set-TextValue "D29" "15.71"
Set-TextValue "E29" "  +0.33%  "
Set-TextValue "D30" "0.0509"
Set-TextValue "E30" "  -0.20%  "
Set-TextValue "D31" "1.19"
Set-TextValue "E31" "  +1.10%  "
Set-TextValue "E32" "  +2.00%  "
Set-TextValue "E33" "  -0.68%  "
Set-TextValue "E34" "  +1.15%  "
Set-TextValue "D35" "1.274.23"
Set-TextValue "E35" "  -1.03%  "
Set-TextValue "D36" "2.44"
Set-TextValue "E36" "  +0.42%  "
Set-TextValue "E37" "  +0.87%  "
Set-TextValue "E38" "  -0.84%  "
Set-TextValue "D39" "0.816"
Set-TextValue "E39" "  -1.26%  "
Set-TextValue "E40" "  -0.18%  "
Set-TextValue "D41" "0.805"
Set-TextValue "E41" "  -0.41%  "
Set-TextValue "E42" "  +0.06%  "
Set-TextValue "D43" "1.783.15"
Set-TextValue "E43" "  -0.64%  "
Set-TextValue "D44" "2.13"
Set-TextValue "E44" "  -4.59%  "
Set-TextValue "D45" "92.83"
Set-TextValue "E45" "  +1.61%  "
Set-TextValue "D46" "61.23"
Set-TextValue "E46" "  +0.15%  "
Set-TextValue "D47" "1.60"
Set-TextValue "E47" "  -0.13%  "
Set-TextValue "E48" "  -1.83%  "
Set-TextValue "E49" "  -0.33%  "
Set-TextValue "D50" "7.57"
Set-TextValue "E50" "  -1.70%  "
Set-TextValue "D51" "0.0966"
Set-TextValue "E51" "  -1.30%  "
